$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
try {
  $excel.ActiveWindow.ScrollRow = 239
  Write-Host "prop scrollrow value:" $excel.ActiveWindow.ScrollRow
} catch {
  Write-Host "err1: $_"
}
try {
    $excel.Goto($ws.Range("A239"), $true)
    Write-Host "goto ok"
} catch {
    Write-Host "err2: $_"
}
